$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.237.40"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "3.169.07"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'604.39"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "'154.10"
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.166.52"
$ws.Range("E8").Value = "  -1.24%  "
$ws.Range("E9").Value = "  +2.73%  "
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("E11").Value = "  -9.17%  "
$ws.Range("D12").Value = "'0.516"
$ws.Range("E12").Value = "  +1.90%  "
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").Value = "'38.38"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "3.689.47"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").Value = "66.209.42"
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "'7.41"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("D18").Value = "3.158.18"
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("D20").Value = "'510.78"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "'15.37"
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("D22").Value = "'0.729"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").Value = "'8.02"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'14.70"
$ws.Range("E24").Value = "  -4.04%  "
$ws.Range("D25").Value = "'84.71"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "'3.00"
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("D28").Value = "'9.19"
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("D29").Value = "'2.40"
$ws.Range("E29").Value = "  +7.04%  "
$ws.Range("D30").Value = "'3.08"
$ws.Range("E30").Value = "  +8.58%  "
$ws.Range("D31").Value = "'7.09"
$ws.Range("E31").Value = "  +4.90%  "
$ws.Range("D32").Value = "'27.96"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -2.02%  "
$ws.Range("E35").Value = "  -1.31%  "
$ws.Range("D36").Value = "'503.34"
$ws.Range("E36").Value = "  +4.80%  "
$ws.Range("D37").Value = "'54.81"
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("D38").Value = "'0.0883"
$ws.Range("E38").Value = "  -2.27%  "
$ws.Range("D39").Value = "'0.0421"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("E40").Value = "  +8.94%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0686"
$ws.Range("E41").Value = "  +7.40%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").Value = "'8.74"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("D43").Value = "'0.297"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").Value = "'2.81"
$ws.Range("E44").Value = "  -4.51%  "
$ws.Range("D45").Value = "'2.46"
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("D46").Value = "2.823.67"
$ws.Range("E46").Value = "  -3.57%  "
$ws.Range("D47").Value = "'28.00"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("E48").Value = "  +3.75%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("D51").Value = "'35.23"
$ws.Range("E51").Value = "  +5.04%  "
